$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 29/30 and 39/40 swap their Coin (B) and Link (C) values as part of
# the underlying source data being re-sorted; apply those first.
$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("B30").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C30").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"

# Updated Volume(1h) (E) values for every data row (2-51). These never look
# like plain numbers (leading/trailing spaces, "%"), so they stay text with
# a plain .Value assignment.
$volumeUpdates = @{
    'E2' = '  +0.23%  '
    'E3' = '  -2.41%  '
    'E4' = '  -1.14%  '
    'E5' = '  +2.60%  '
    'E6' = '  -0.44%  '
    'E7' = '  -2.17%  '
    'E8' = '  -4.50%  '
    'E9' = '  -4.08%  '
    'E10' = '  -5.65%  '
    'E11' = '  -5.22%  '
    'E12' = '  -1.05%  '
    'E13' = '  -3.92%  '
    'E14' = '  -5.69%  '
    'E15' = '  -2.45%  '
    'E16' = '  -5.31%  '
    'E17' = '  -6.73%  '
    'E18' = '  -1.33%  '
    'E19' = '  -0.30%  '
    'E20' = '  -4.90%  '
    'E21' = '  -6.07%  '
    'E22' = '  -7.48%  '
    'E23' = '  -2.56%  '
    'E24' = '  -0.32%  '
    'E25' = '  -0.63%  '
    'E26' = '  -12.89%  '
    'E27' = '  -0.24%  '
    'E28' = '  -7.36%  '
    'E29' = '  -1.58%  '
    'E30' = '  -2.63%  '
    'E31' = '  -4.25%  '
    'E32' = '  -3.48%  '
    'E33' = '  -12.23%  '
    'E34' = '  -3.23%  '
    'E35' = '  -5.32%  '
    'E36' = '  -9.24%  '
    'E37' = '  -5.33%  '
    'E38' = '  -0.12%  '
    'E39' = '  -6.06%  '
    'E40' = '  -8.52%  '
    'E41' = '  -5.93%  '
    'E42' = '  -8.90%  '
    'E43' = '  -0.20%  '
    'E44' = '  -6.73%  '
    'E45' = '  +0.47%  '
    'E46' = '  -7.37%  '
    'E47' = '  -6.96%  '
    'E48' = '  -2.73%  '
    'E49' = '  -6.44%  '
    'E50' = '  -3.97%  '
    'E51' = '  -1.61%  '
}

foreach ($addr in $volumeUpdates.Keys) {
    $ws.Range($addr).Value = $volumeUpdates[$addr]
}

# Updated Price (D) values for every data row (2-51). Many of these look
# like plain decimal numbers to Excel's auto-detection (e.g. "0.9972"), but
# the source column is plain text, so force the whole column to Text format
# before writing, then restore General/Normal so no stray number format is
# left behind on the cells.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$priceUpdates = @{
    'D2' = '24.822.31'
    'D3' = '1.659.95'
    'D4' = '0.9972'
    'D5' = '319.19'
    'D6' = '0.9982'
    'D7' = '0.3640'
    'D8' = '46.86'
    'D9' = '0.3277'
    'D10' = '1.141'
    'D11' = '0.07085'
    'D12' = '0.9933'
    'D13' = '6.046'
    'D14' = '19.73'
    'D15' = '1.658.68'
    'D16' = '6.633'
    'D17' = '0.00001049'
    'D18' = '0.06625'
    'D19' = '0.9986'
    'D20' = '79.36'
    'D21' = '5.942'
    'D22' = '15.86'
    'D23' = '12.64'
    'D24' = '24.699.20'
    'D25' = '2.437'
    'D26' = '2.405'
    'D27' = '148.59'
    'D28' = '18.72'
    'D29' = '1.222'
    'D30' = '1.838.58'
    'D31' = '125.94'
    'D32' = '4.084'
    'D33' = '5.897'
    'D34' = '0.08439'
    'D35' = '1.671'
    'D36' = '12.36'
    'D37' = '5.252'
    'D38' = '1.269'
    'D39' = '0.02250'
    'D40' = '0.06051'
    'D41' = '0.2085'
    'D42' = '8.248'
    'D43' = '0.9999'
    'D44' = '0.5960'
    'D45' = '3.838'
    'D46' = '12.74'
    'D47' = '0.5672'
    'D48' = '125.37'
    'D49' = '1.962'
    'D50' = '0.06973'
    'D51' = '1.197'
}

foreach ($addr in $priceUpdates.Keys) {
    $ws.Range($addr).Value = $priceUpdates[$addr]
}

$priceRange.NumberFormat = "General"
$priceRange.Style = "Normal"
